$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the project hours value for week 15 (row 15, column F)
$ws.Range("F15").Value = 2

# Move the active selection to E15 (matches the recorded cursor position change)
$ws.Range("E15").Select()

$wb.Save()
